# Applies the "Updated cryptos list" data refresh described in the commit:
# refreshed Price (D) / Volume(1h) (E) figures for rows 2-51, plus the
# Hedera <-> PancakeSwap row swap (rows 30-31, columns B/C/D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several Price cells hold numeric-looking text (e.g. "0.485", "26.765.75")
# that must stay as plain text. Forcing a text NumberFormat before the
# assignment stops Excel from auto-converting them to real numbers; the
# style is then restored so cell formatting is left untouched.
function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.788.40"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.539.43"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "205.87"
$ws.Range("E5").Value = "  -0.28%  "
Set-TextValue "D6" "0.486"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue "D8" "0.247"
$ws.Range("E8").Value = "  -0.30%  "
Set-TextValue "D9" "21.29"
$ws.Range("E9").Value = "  -2.83%  "
Set-TextValue "D10" "0.0581"
$ws.Range("E10").Value = "  -0.47%  "
Set-TextValue "D11" "0.0854"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "1.757.96"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "1.540.33"
$ws.Range("E13").Value = "  -1.82%  "
Set-TextValue "D14" "3.68"
$ws.Range("E14").Value = "  -1.28%  "
Set-TextValue "D15" "0.509"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "26.789.81"
$ws.Range("E16").Value = "  +0.01%  "
Set-TextValue "D17" "61.10"
$ws.Range("E17").Value = "  -0.40%  "
Set-TextValue "D18" "213.68"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").Value = "0.0₃0683"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -2.46%  "
Set-TextValue "D23" "9.17"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("E24").Value = "  -3.48%  "
Set-TextValue "D25" "151.79"
$ws.Range("E25").Value = "  -0.43%  "
Set-TextValue "D26" "6.58"
$ws.Range("E26").Value = "  -2.25%  "
Set-TextValue "D27" "14.79"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.0459"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "1.10"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "1.366.69"
$ws.Range("E33").Value = "  -1.99%  "
Set-TextValue "D34" "2.94"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  -1.46%  "
Set-TextValue "D36" "0.960"
$ws.Range("E36").Value = "  +3.03%  "
Set-TextValue "D37" "2.28"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +0.88%  "
Set-TextValue "D39" "0.521"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  +8.06%  "
Set-TextValue "D41" "0.805"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  +1.02%  "
Set-TextValue "D44" "62.95"
$ws.Range("E44").Value = "  -0.43%  "
Set-TextValue "D45" "1.74"
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("D46").Value = "1.672.35"
$ws.Range("E46").Value = "  -1.69%  "
Set-TextValue "D47" "84.36"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("D49").Value = "0.0₇0969"
$ws.Range("E49").Value = "  -1.86%  "
Set-TextValue "D50" "0.0943"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -0.04%  "
